# Update "想去人数" (want-to-go headcount) values across the four sheets:
#   Sheet 1 "展览"   (Exhibitions)
#   Sheet 2 "演出"   (Performances)
#   Sheet 3 "本地生活" (Local life)
#   Sheet 4 "全部类型" (All types, a combined roll-up of sheets 1-3)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value  = 1677
$ws1.Range("F3").Value  = 9587
$ws1.Range("F10").Value = 1445
$ws1.Range("F13").Value = 1513
$ws1.Range("F17").Value = 158
$ws1.Range("F39").Value = 339

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F24").Value = 36
$ws2.Range("F39").Value = 35

$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F6").Value  = 2426
$ws3.Range("F7").Value  = 3763
$ws3.Range("F10").Value = 138
$ws3.Range("F11").Value = 126

$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value  = 1677
$ws4.Range("F5").Value  = 9587
$ws4.Range("F6").Value  = 3763
$ws4.Range("F8").Value  = 138
$ws4.Range("F9").Value  = 138
$ws4.Range("F16").Value = 1445
$ws4.Range("F18").Value = 126
$ws4.Range("F19").Value = 1513
$ws4.Range("F23").Value = 158
$ws4.Range("F36").Value = 36
$ws4.Range("F41").Value = 339
$ws4.Range("F50").Value = 35
